$d = $word.ActiveDocument

# The footer block at the end of the document (an empty spacer
# paragraph, "Ver no Jupiter Salvar em pdf Salvar em docx", and the
# "© 2020 ... Creative Commons Attribution" copyright line) is being
# removed, leaving the "LOT2007: ..." requisitos paragraph followed
# directly by the existing empty paragraph that precedes the trailing
# page break.

$jupiterIdx = -1
$copyrightIdx = -1

$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    $t = $p.Range.Text
    if ($t -like "*Ver no Jupiter*") {
        $jupiterIdx = $idx
    }
    if ($t -like "*Creative Commons Attribution*") {
        $copyrightIdx = $idx
    }
}

if ($jupiterIdx -gt 0 -and $copyrightIdx -ge $jupiterIdx) {
    # Include the blank paragraph immediately before the "Ver no
    # Jupiter..." line so all three paragraphs collapse away together.
    $startPara = $d.Paragraphs.Item($jupiterIdx - 1)
    $endPara = $d.Paragraphs.Item($copyrightIdx)

    $deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $deleteRange.Delete()
}
